$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "humanitarian_scope_pct" column (F) is empty for every data row.
# Delete it entirely, which shifts humanitarian_scope_desc_eng (G) into F
# and humanitarian_scope_desc_fr (H) into G.
$ws.Columns("F").Delete()
